$d = $word.ActiveDocument

# The document currently ends with a single, empty trailing paragraph
# (plain <w:p/>, no pPr). Per the target diff we need that paragraph to
# gain a pPr/rPr/rFonts hint="eastAsia", and two new content paragraphs
# need to be appended right after it (the second of which becomes the
# document's new final paragraph).
#
# Word will not let us delete the very last paragraph mark in a
# document, so instead of trying to edit the existing trailing
# paragraph in place, we insert all three target paragraphs (the
# pPr-carrying empty one + the two content ones) *before* the existing
# trailing paragraph, then merge the now-redundant original trailing
# paragraph away (it is no longer the document's final paragraph mark,
# so it can be deleted/merged normally).

$lastPara = $d.Paragraphs.Last
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)  # wdCollapseStart

$newContentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">When the K-D tree is constructed, it can be very </w:t></w:r><w:r><w:t xml:space="preserve">useful for problem like nearest neighbor </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">search and range search. When in range search, start at the rood, if the range covers(surpass) the hyperplane , both the branch must be search. </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>f not, only one branch is need to be searched.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">When is nearest neighbor search , it first search till the leaf, than goes back and see if the current min distance covers the hyperplane or not. </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>f not, goes back one node, if yes, it recursively into the other branch.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($newContentXml)

# Merge the now-redundant original trailing empty paragraph into the
# paragraph immediately before it (our new last content paragraph),
# by deleting that paragraph's own ending mark. Word keeps the
# *surviving* (later) paragraph mark's formatting, and since both the
# original trailing paragraph and our new last content paragraph use
# default paragraph formatting, this merge is formatting-neutral while
# removing the extra empty paragraph.
$paraCount = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($paraCount - 1)
$markRange = $d.Range($secondToLast.Range.End - 1, $secondToLast.Range.End)
$markRange.Delete()

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
